# Auto-generated edit script: updates profit-calculation columns (H-N)
# across multiple job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to reflect refreshed market-board pricing data.

$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC!row46 (Leve Item ID 4584)
$ws_ALC.Range("H46").Value = 1500
$ws_ALC.Range("J46").Value = 1500
$ws_ALC.Range("L46").Value = 4500
$ws_ALC.Range("N46").Value = -4738

# ALC!row60 (Leve Item ID 4584)
$ws_ALC.Range("H60").Value = 1500
$ws_ALC.Range("J60").Value = 1500
$ws_ALC.Range("L60").Value = 4500
$ws_ALC.Range("N60").Value = -5468

# ALC!row70 (Leve Item ID 12604)
$ws_ALC.Range("H70").Value = 2624.75
$ws_ALC.Range("I70").Value = 2299
$ws_ALC.Range("J70").Value = 2733.3333
$ws_ALC.Range("K70").Value = 6897
$ws_ALC.Range("L70").Value = 8199.999899999999
$ws_ALC.Range("M70").Value = -6627
$ws_ALC.Range("N70").Value = -8739.999899999999

# ALC!row73 (Leve Item ID 12604)
$ws_ALC.Range("H73").Value = 2624.75
$ws_ALC.Range("I73").Value = 2299
$ws_ALC.Range("J73").Value = 2733.3333
$ws_ALC.Range("K73").Value = 6897
$ws_ALC.Range("L73").Value = 8199.999899999999
$ws_ALC.Range("M73").Value = -5961
$ws_ALC.Range("N73").Value = -10071.9999

# ALC!row99 (Leve Item ID 19883)
$ws_ALC.Range("H99").Value = 342.33334
$ws_ALC.Range("J99").Value = 456.5
$ws_ALC.Range("L99").Value = 1369.5
$ws_ALC.Range("N99").Value = -4365.5

# ALC!row125 (Leve Item ID 36228)
$ws_ALC.Range("H125").Value = 125005050
$ws_ALC.Range("I125").Value = 142861630
$ws_ALC.Range("K125").Value = 1285754670
$ws_ALC.Range("M125").Value = -1285752210

# ALC!row137 (Leve Item ID 44013)
$ws_ALC.Range("H137").Value = 1143.8572
$ws_ALC.Range("I137").Value = 1143.8572
$ws_ALC.Range("J137").Value = 0
$ws_ALC.Range("K137").Value = 3431.5716
$ws_ALC.Range("L137").Value = 0
$ws_ALC.Range("M137").Value = -881.5715999999998
$ws_ALC.Range("N137").ClearContents()

# ALC!row138 (Leve Item ID 44169)
$ws_ALC.Range("H138").Value = 3927.2307
$ws_ALC.Range("I138").Value = 1853.1875
$ws_ALC.Range("J138").Value = 5370.0435
$ws_ALC.Range("K138").Value = 5559.5625
$ws_ALC.Range("L138").Value = 16110.1305
$ws_ALC.Range("M138").Value = -419.5625
$ws_ALC.Range("N138").Value = -26390.1305

# ARM!row2 (Leve Item ID 27713)
$ws_ARM.Range("H2").Value = 1861
$ws_ARM.Range("I2").Value = 1861
$ws_ARM.Range("K2").Value = 1861
$ws_ARM.Range("M2").Value = -1748

# ARM!row45 (Leve Item ID 27714)
$ws_ARM.Range("H45").Value = 84465.664
$ws_ARM.Range("I45").Value = 143984.14
$ws_ARM.Range("K45").Value = 143984.14
$ws_ARM.Range("M45").Value = -143607.14

# ARM!row63 (Leve Item ID 12528)
$ws_ARM.Range("H63").Value = 5358
$ws_ARM.Range("I63").Value = 5271.1665
$ws_ARM.Range("J63").Value = 6400
$ws_ARM.Range("K63").Value = 5271.1665
$ws_ARM.Range("L63").Value = 6400
$ws_ARM.Range("M63").Value = -4585.1665
$ws_ARM.Range("N63").Value = -7772

# ARM!row66 (Leve Item ID 12528)
$ws_ARM.Range("H66").Value = 5358
$ws_ARM.Range("I66").Value = 5271.1665
$ws_ARM.Range("J66").Value = 6400
$ws_ARM.Range("K66").Value = 26355.8325
$ws_ARM.Range("L66").Value = 32000
$ws_ARM.Range("M66").Value = -22923.8325
$ws_ARM.Range("N66").Value = -38864

# ARM!row97 (Leve Item ID 19941)
$ws_ARM.Range("H97").Value = 1999.5
$ws_ARM.Range("I97").Value = 2232.6667
$ws_ARM.Range("J97").Value = 1766.3334
$ws_ARM.Range("K97").Value = 2232.6667
$ws_ARM.Range("L97").Value = 1766.3334
$ws_ARM.Range("M97").Value = -1736.6667
$ws_ARM.Range("N97").Value = -2758.3334

# ARM!row116 (Leve Item ID 27713)
$ws_ARM.Range("H116").Value = 1861
$ws_ARM.Range("I116").Value = 1861
$ws_ARM.Range("K116").Value = 1861
$ws_ARM.Range("M116").Value = 433

# BSM!row3 (Leve Item ID 27713)
$ws_BSM.Range("H3").Value = 1861
$ws_BSM.Range("I3").Value = 1861
$ws_BSM.Range("K3").Value = 1861
$ws_BSM.Range("M3").Value = -1747

# BSM!row110 (Leve Item ID 25790)
$ws_BSM.Range("H110").Value = 50000
$ws_BSM.Range("J110").Value = 50000
$ws_BSM.Range("L110").Value = 50000
$ws_BSM.Range("N110").Value = -58180

# BSM!row134 (Leve Item ID 43998)
$ws_BSM.Range("H134").Value = 3469.25
$ws_BSM.Range("I134").Value = 3469.25
$ws_BSM.Range("K134").Value = 10407.75
$ws_BSM.Range("M134").Value = -7872.75

# CRP!row31 (Leve Item ID 44023)
$ws_CRP.Range("H31").Value = 1288.12
$ws_CRP.Range("I31").Value = 1015.9474
$ws_CRP.Range("J31").Value = 2150
$ws_CRP.Range("K31").Value = 1015.9474
$ws_CRP.Range("L31").Value = 2150
$ws_CRP.Range("M31").Value = -720.9474
$ws_CRP.Range("N31").Value = -2740

# CRP!row34 (Leve Item ID 44023)
$ws_CRP.Range("H34").Value = 1288.12
$ws_CRP.Range("I34").Value = 1015.9474
$ws_CRP.Range("J34").Value = 2150
$ws_CRP.Range("K34").Value = 1015.9474
$ws_CRP.Range("L34").Value = 2150
$ws_CRP.Range("M34").Value = -813.9474
$ws_CRP.Range("N34").Value = -2554

# CRP!row43 (Leve Item ID 18504)
$ws_CRP.Range("H43").Value = 25000
$ws_CRP.Range("J43").Value = 25000
$ws_CRP.Range("L43").Value = 25000
$ws_CRP.Range("N43").Value = -25368

# CRP!row101 (Leve Item ID 18504)
$ws_CRP.Range("H101").Value = 25000
$ws_CRP.Range("J101").Value = 25000
$ws_CRP.Range("L101").Value = 25000
$ws_CRP.Range("N101").Value = -31490

# CRP!row105 (Leve Item ID 19928)
$ws_CRP.Range("H105").Value = 3490.4167
$ws_CRP.Range("J105").Value = 4328.5713
$ws_CRP.Range("L105").Value = 4328.5713
$ws_CRP.Range("N105").Value = -7822.5713

# CUL!row6 (Leve Item ID 4639)
$ws_CUL.Range("H6").Value = 12945.75
$ws_CUL.Range("I6").Value = 14581.571
$ws_CUL.Range("J6").Value = 1495
$ws_CUL.Range("K6").Value = 43744.713
$ws_CUL.Range("L6").Value = 4485
$ws_CUL.Range("M6").Value = -43631.713
$ws_CUL.Range("N6").Value = -4711

# CUL!row39 (Leve Item ID 4712)
$ws_CUL.Range("H39").Value = 7500
$ws_CUL.Range("J39").Value = 7500
$ws_CUL.Range("L39").Value = 22500
$ws_CUL.Range("N39").Value = -23088

# CUL!row55 (Leve Item ID 4733)
$ws_CUL.Range("H55").Value = 499
$ws_CUL.Range("I55").Value = 499
$ws_CUL.Range("K55").Value = 1497
$ws_CUL.Range("M55").Value = -1320

# CUL!row107 (Leve Item ID 27838)
$ws_CUL.Range("H107").Value = 0
$ws_CUL.Range("I107").Value = 0
$ws_CUL.Range("J107").Value = 0
$ws_CUL.Range("K107").Value = 0
$ws_CUL.Range("L107").Value = 0
$ws_CUL.Range("M107").ClearContents()
$ws_CUL.Range("N107").ClearContents()

# CUL!row112 (Leve Item ID 27855)
$ws_CUL.Range("H112").Value = 8907.25
$ws_CUL.Range("J112").Value = 14332.333
$ws_CUL.Range("L112").Value = 42996.999
$ws_CUL.Range("N112").Value = -45212.999

# CUL!row125 (Leve Item ID 36043)
$ws_CUL.Range("H125").Value = 0
$ws_CUL.Range("I125").Value = 0
$ws_CUL.Range("K125").Value = 0
$ws_CUL.Range("M125").ClearContents()

# CUL!row131 (Leve Item ID 36060)
$ws_CUL.Range("H131").Value = 418636
$ws_CUL.Range("I131").Value = 907.1667
$ws_CUL.Range("J131").Value = 557878.9399999999
$ws_CUL.Range("K131").Value = 2721.5001
$ws_CUL.Range("L131").Value = 1673636.82
$ws_CUL.Range("M131").Value = 2318.4999
$ws_CUL.Range("N131").Value = -1683716.82

# CUL!row132 (Leve Item ID 43972)
$ws_CUL.Range("H132").Value = 700
$ws_CUL.Range("J132").Value = 700
$ws_CUL.Range("L132").Value = 6300
$ws_CUL.Range("N132").Value = -11360

# GSM!row132 (Leve Item ID 44008)
$ws_GSM.Range("H132").Value = 1807.6666
$ws_GSM.Range("I132").Value = 1807.6666
$ws_GSM.Range("K132").Value = 5422.9998
$ws_GSM.Range("M132").Value = -2892.9998

# LTW!row7 (Leve Item ID 36249)
$ws_LTW.Range("H7").Value = 6754.409
$ws_LTW.Range("I7").Value = 3651.75
$ws_LTW.Range("J7").Value = 7443.8887
$ws_LTW.Range("K7").Value = 3651.75
$ws_LTW.Range("L7").Value = 7443.8887
$ws_LTW.Range("M7").Value = -3539.75
$ws_LTW.Range("N7").Value = -7667.8887

# LTW!row40 (Leve Item ID 36248)
$ws_LTW.Range("H40").Value = 3824.25
$ws_LTW.Range("I40").Value = 3265.6667
$ws_LTW.Range("K40").Value = 3265.6667
$ws_LTW.Range("M40").Value = -3129.6667

# LTW!row112 (Leve Item ID 25846)
$ws_LTW.Range("H112").Value = 52989.5
$ws_LTW.Range("J112").Value = 52989.5
$ws_LTW.Range("L112").Value = 52989.5
$ws_LTW.Range("N112").Value = -55943.5

# LTW!row126 (Leve Item ID 36249)
$ws_LTW.Range("H126").Value = 6754.409
$ws_LTW.Range("I126").Value = 3651.75
$ws_LTW.Range("J126").Value = 7443.8887
$ws_LTW.Range("K126").Value = 10955.25
$ws_LTW.Range("L126").Value = 22331.6661
$ws_LTW.Range("M126").Value = -8485.25
$ws_LTW.Range("N126").Value = -27271.6661

# LTW!row132 (Leve Item ID 44058)
$ws_LTW.Range("H132").Value = 15922.5
$ws_LTW.Range("I132").Value = 26842
$ws_LTW.Range("J132").Value = 5003
$ws_LTW.Range("K132").Value = 80526
$ws_LTW.Range("L132").Value = 15009
$ws_LTW.Range("M132").Value = -77996
$ws_LTW.Range("N132").Value = -20069

# WVR!row81 (Leve Item ID 12596)
$ws_WVR.Range("H81").Value = 2000575
$ws_WVR.Range("I81").Value = 624.6667
$ws_WVR.Range("J81").Value = 5000500.5
$ws_WVR.Range("K81").Value = 1249.3334
$ws_WVR.Range("L81").Value = 10001001
$ws_WVR.Range("M81").Value = -188.3334
$ws_WVR.Range("N81").Value = -10003123

# WVR!row84 (Leve Item ID 12596)
$ws_WVR.Range("H84").Value = 2000575
$ws_WVR.Range("I84").Value = 624.6667
$ws_WVR.Range("J84").Value = 5000500.5
$ws_WVR.Range("K84").Value = 6246.666999999999
$ws_WVR.Range("L84").Value = 50005005
$ws_WVR.Range("M84").Value = -942.6669999999995
$ws_WVR.Range("N84").Value = -50015613

# WVR!row126 (Leve Item ID 36210)
$ws_WVR.Range("H126").Value = 2674.75
$ws_WVR.Range("I126").Value = 2674.75
$ws_WVR.Range("K126").Value = 8024.25
$ws_WVR.Range("M126").Value = -5554.25

# WVR!row132 (Leve Item ID 44029)
$ws_WVR.Range("H132").Value = 3326.6667
$ws_WVR.Range("I132").Value = 5047.6
$ws_WVR.Range("J132").Value = 1175.5
$ws_WVR.Range("K132").Value = 15142.8
$ws_WVR.Range("L132").Value = 3526.5
$ws_WVR.Range("M132").Value = -12612.8
$ws_WVR.Range("N132").Value = -8586.5

# WVR!row136 (Leve Item ID 44031)
$ws_WVR.Range("H136").Value = 1889.9546
$ws_WVR.Range("I136").Value = 1854.3684
$ws_WVR.Range("K136").Value = 5563.1052
$ws_WVR.Range("M136").Value = -3013.1052

